$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("studios")
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("B10").Select() | Out-Null
Write-Host "done"
